$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.45"
$ws.Range("E2").Value = "'-0.10%"
$ws.Range("D3").Value = "'26.91"
$ws.Range("E3").Value = "'-0.05%"
$ws.Range("D4").Value = "'4.805"
$ws.Range("E4").Value = "'1.89%"
$ws.Range("D5").Value = "'0.05946"
$ws.Range("E5").Value = "'0.25%"
$ws.Range("D6").Value = "'6.606"
$ws.Range("E6").Value = "'-0.82%"
$ws.Range("D7").Value = "'0.8505"
$ws.Range("E7").Value = "'-1.92%"
$ws.Range("D8").Value = "'0.9247"
$ws.Range("E8").Value = "'-1.34%"
$ws.Range("D9").Value = "'0.1379"
$ws.Range("E9").Value = "'-1.39%"
$ws.Range("D10").Value = "'0.04195"
$ws.Range("E10").Value = "'12.09%"
$ws.Range("D11").Value = "'0.07010"
$ws.Range("E11").Value = "'-1.51%"
$ws.Range("E12").Value = "'-3.45%"
$ws.Range("D13").Value = "'0.09092"
$ws.Range("E13").Value = "'-1.66%"
$ws.Range("D14").Value = "'0.001527"
$ws.Range("E14").Value = "'-0.48%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006088"
$ws.Range("E15").Value = "'1.46%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.471"
$ws.Range("E16").Value = "'-0.58%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.161"
$ws.Range("E17").Value = "'-1.32%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.199"
$ws.Range("E18").Value = "'-1.85%"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "'0.01029"
$ws.Range("E19").Value = "'1,599.43%"
$ws.Range("D20").Value = "'0.3039"
$ws.Range("E20").Value = "'-2.79%"
$ws.Range("E21").Value = "'0.28%"
$ws.Range("D22").Value = "'3.917"
$ws.Range("E22").Value = "'2.93%"
$ws.Range("D23").Value = "'0.04233"
$ws.Range("E23").Value = "'0.43%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'-0.15%"
$ws.Range("D25").Value = "'0.003619"
$ws.Range("E25").Value = "'-15.56%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.07%"
$ws.Range("E27").Value = "'1.97%"
$ws.Range("E40").Value = "'-1.04%"
$ws.Range("D41").Value = "'0.006278"
$ws.Range("E41").Value = "'1.52%"
$ws.Range("D42").Value = "'0.1097"
$ws.Range("E42").Value = "'-0.34%"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'-2.23%"
$ws.Range("D44").Value = "'0.01413"
$ws.Range("E44").Value = "'26.69%"
$ws.Range("D45").Value = "'0.00005343"
$ws.Range("E45").Value = "'-2.76%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("D47").Value = "'0.04398"
$ws.Range("E47").Value = "'-50.26%"
$ws.Range("E48").Value = "'9,870.25%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.07%"
